$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing Writing marks (I11, I12)
$ws.Range("I11").Value = 6
$ws.Range("I12").Value = 4

# Row 14: new IELTS11_Test2 entry (Listening only) - set first so its
# shared string lands before the "a Practice Test" one (matches author's order)
$ws.Range("C14").Value = [DateTime]"2024-06-04"
$ws.Range("D14").Value = "IELTS11_Test2"
$ws.Range("E14").Value = 30
$ws.Range("F14").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Lis_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'

# Row 13: new Practice Test entry (Reading only)
$ws.Range("C13").Value = [DateTime]"2024-05-31"
$ws.Range("D13").Value = "a Practice Test"
$ws.Range("G13").Value = 25
$ws.Range("H13").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'

# Update selection to match final state
$ws.Range("H14").Select()
